$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.338.02"
$ws.Range("E2").Value = "  -3.93%  "

$ws.Range("D3").Value = "2.615.42"
$ws.Range("E3").Value = "  -3.32%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.74%  "

$ws.Range("E10").Value = "  -3.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "3.074.52"
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("D14").Value = "58.293.86"
$ws.Range("E14").Value = "  -3.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.10%  "

$ws.Range("D17").Value = "2.598.74"
$ws.Range("E17").Value = "  -4.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "336.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.59%  "

$ws.Range("E19").Value = "  -3.10%  "

$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.414"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "

$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("D28").Value = "0.0₃0788"
$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.96%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.19%  "

$ws.Range("E35").Value = "  -4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.877"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0966"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "268.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0530"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "

$ws.Range("D48").Value = "2.020.67"
$ws.Range("E48").Value = "  -5.10%  "

$ws.Range("E49").Value = "  -2.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.88%  "
